$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts existing rows 21-42 down to 22-43)
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new weekly record
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C21").Value = "Los Lagos"
$ws.Range("D21").Value = 45174
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 100112013
$ws.Range("G21").Value = "Alcachofa"
$ws.Range("H21").Value = "Argentina(o)"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 15000
$ws.Range("N21").Value = '$/caja 50 unidades'
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 300
$ws.Range("Q21").Value = 50
$ws.Range("R21").Value = "Hortaliza"
